{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the last \"group\" list paragraph - the one ending with\n// \"Md Nazmum Hasan Nafees.\" - and append the two new group entries\n// right after it (still before the trailing blank paragraph).\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text || \"\";\n  if (text.indexOf(\"Nazmum Hasan Nafees\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate anchor paragraph for insertion.\");\n}\n\nconst firstNew = anchor.insertParagraph(\n  \"Jamaal Nnamdi, Mahtab Askarzadeh, Maxwell Jones.\",\n  Word.InsertLocation.after\n);\nconst secondNew = firstNew.insertParagraph(\n  \"Arvin Armand, Arvin Salehi.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the last existing group entry - the paragraph ending with\n# \"Md Nazmum Hasan Nafees.\" - so the two new group lines get appended\n# right after it (still before the trailing blank paragraph).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Nazmum Hasan Nafees*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not locate anchor paragraph for insertion.\"\n}\n\n$target.Range.InsertParagraphAfter()\n$firstNew = $target.Next()\n$firstNew.Range.Text = \"Jamaal Nnamdi, Mahtab Askarzadeh, Maxwell Jones.\"\n\n$firstNew.Range.InsertParagraphAfter()\n$secondNew = $firstNew.Next()\n$secondNew.Range.Text = \"Arvin Armand, Arvin Salehi.\"\n"}
